$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column A (temperature) to include units
$ws.Range("A1").Value = "T, K"

# Update bias-field / time series values in column B (rows 2-17)
$bValues = @(10, 15, 18, 21, 24, 27, 30, 33, 36, 39, 42, 45, 48, 51, 54, 57)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Update C9 value
$ws.Range("C9").Value = 10

# Update the active selection shown when the workbook is opened
$ws.Range("J10").Select()
